$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 29 (S.No 18, CustomerMappingDriver Class section) ---
# Score for "Getter methods" / findNoOfCustomers-style row changed from 8 to 15
$ws.Range("E29").Value = 15
# Grading comment updated to reflect a new deduction reason
$ws.Range("F29").Value = "(-1) for using == when comparing string instead of .equals() method"

# --- Row 30 (S.No 19) ---
$ws.Range("F30").Value = "(-4) for incorrect output due to ArrayIndexOutOfBoundsException"

# --- Row 37 (compilation errors deduction row) ---
$ws.Range("F37").Value = "(-2.5) for getting ArrayIndexOutOfBoundsException"

# The dependent SUM totals (E31, E38) recalculate automatically.

# Reflect where the grader ended up working when they saved the file
$ws.Activate()
$ws.Range("F29").Select()
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 1
